# Applies the cryptos.xlsx price/volume refresh (and the three-row coin
# re-ranking at rows 44-51) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text that frequently *looks* numeric (e.g. "6.51",
# "40.172.56" as a thousands-grouped price). Assigning such a string straight
# to .Value lets Excel reinterpret it as a real number, which would both
# change the stored cell type and could mangle multi-dot price strings.
# Prefixing with an apostrophe forces text entry; re-applying the "Normal"
# style afterwards drops the implicit @ (text) number-format Excel tacks on,
# so the cell keeps the workbook's original (default) style.
function Set-TextCell($ws, $ref, $value) {
    $c = $ws.Range($ref)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

# --- Coin / link text cells (row re-ranking, rows 44-51) ---
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"

# --- Price cells (column D) ---
Set-TextCell $ws "D2" "40.172.56"
Set-TextCell $ws "D3" "2.229.03"
Set-TextCell $ws "D5" "292.09"
Set-TextCell $ws "D6" "88.26"
Set-TextCell $ws "D9" "0.475"
Set-TextCell $ws "D10" "30.60"
Set-TextCell $ws "D13" "6.51"
Set-TextCell $ws "D14" "2.573.03"
Set-TextCell $ws "D15" "14.04"
Set-TextCell $ws "D16" "2.232.16"
Set-TextCell $ws "D17" "0.730"
Set-TextCell $ws "D18" "40.128.43"
Set-TextCell $ws "D19" "11.54"
Set-TextCell $ws "D22" "65.78"
Set-TextCell $ws "D23" "236.68"
Set-TextCell $ws "D26" "1.84"
Set-TextCell $ws "D27" "22.77"
Set-TextCell $ws "D30" "155.90"
Set-TextCell $ws "D31" "32.01"
Set-TextCell $ws "D33" "4.97"
Set-TextCell $ws "D34" "0.0719"
Set-TextCell $ws "D35" "2.34"
Set-TextCell $ws "D36" "2.89"
Set-TextCell $ws "D38" "15.71"
Set-TextCell $ws "D39" "0.0987"
Set-TextCell $ws "D40" "1.71"
Set-TextCell $ws "D41" "2.111.28"
Set-TextCell $ws "D43" "2.15"
Set-TextCell $ws "D44" "0.0269"
Set-TextCell $ws "D45" "17.95"
Set-TextCell $ws "D46" "9.99"
Set-TextCell $ws "D48" "2.437.88"
Set-TextCell $ws "D49" "1.46"
Set-TextCell $ws "D50" "89.32"
Set-TextCell $ws "D51" "69.56"

# --- Volume(1h) percentage cells (column E) ---
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("E12").Value = "  +3.28%  "
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  +7.39%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  -5.04%  "
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  +6.92%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -5.97%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +7.90%  "
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  +10.01%  "
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("E51").Value = "  -2.25%  "
